# Quarterly indexing esoteric bug-fix operation
#
# Column A holds the "as-of" date for each forecast row. Each date was
# incorrectly stamped as the 1st of its quarter-start month; the fix
# re-stamps every date in column A to the 15th of the month that follows
# the originally recorded month (e.g. 2008-01-01 -> 2008-02-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 73
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $rawValue = $cell.Value2

    if ($rawValue -eq $null -or $rawValue -eq "") {
        continue
    }

    $oldSerial = [double]$rawValue
    $oldDate = [DateTime]::FromOADate($oldSerial)

    $shifted = $oldDate.AddMonths(1)
    $newDate = Get-Date -Year $shifted.Year -Month $shifted.Month -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value = $newDate.ToOADate()
}
